# Update student_information.xlsx: replace existing row 2 data and
# append two additional student rows (rows 3 and 4).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: A Full Name, B Email, C Home Number, D Cell Number, E Campus,
#          F Grade, G Gender, H Date of Birth, I Address,
#          J Registration Date, K State, L Remark
$data = @(
    @{ Row=2; Values=@("Janet Scott", "deanna94@example.org", "557.987.3077x47839", "6462397282", "Henrymouth", 3, "Female", "2013-04-17", "1599 Amanda Plaza Suite 627, East Victoria, PW 17023", "2024-02-02", "South Dakota", "None") },
    @{ Row=3; Values=@("Lynn Flores", "charles46@example.org", "+1-334-410-1697x5063", "+1-813-622-9373", "Jordanville", 12, "Male", "2006-07-02", "62931 Chelsea Shore Suite 104, East George, IL 63484", "2024-01-31", "North Carolina", "None") },
    @{ Row=4; Values=@("Frank Castro", "hernandezcrystal@example.com", "(339)858-7240x786", "810-714-0034", "West Jeffrey", 7, "Female", "2014-03-23", "93634 James Lane, Hansenchester, MS 17361", "2024-02-01", "Wyoming", "None") }
)

# Columns that hold numbers/dates as plain text in the source data and
# must not be auto-converted to numeric or date cell values.
$textColumns = @(4, 8, 10)

foreach ($entry in $data) {
    $r = $entry.Row
    $values = $entry.Values
    for ($i = 0; $i -lt $values.Length; $i++) {
        $col = $i + 1
        $cell = $ws.Cells.Item($r, $col)
        if ($textColumns -contains $col) {
            # Force text storage so strings like "6462397282" or
            # "2013-04-17" are kept literally instead of being parsed
            # into a number / date serial value.
            $cell.NumberFormat = "@"
            $cell.Value = $values[$i]
            $cell.Style = "Normal"
        } else {
            $cell.Value = $values[$i]
        }
    }
}
